# email update for request reroute
# Adds a new "rerouted shipment" data row's worth of values to the
# Input sheet: row 9 (OrderId/pickUpDate/Amount/Tracking#/WayBill) gets a
# fresh UPS record, and row 2 gets a fresh PitneyBowes (PBID) record -
# these simulate a freshly appended QuickQuote data row used by the
# RequestReroute QA scenario.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

function Set-RerouteCell($addr, $text) {
    $cell = $ws.Range($addr)
    # Leading single-quote forces literal text (no date/number/currency
    # auto-conversion) while the value is applied; re-asserting the
    # "Normal" style immediately afterward drops the transient
    # quote-prefix marker so the cell settles back to a plain General
    # text cell identical in kind to its neighbours.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
    # Re-apply the same "new row" banding used by the rest of the
    # table: solid white fill plus a thin top/bottom rule.
    $cell.Interior.ColorIndex = 2
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
}

# Row 9 - new UPS shipment (OrderId, pickUpDate, Amount, Tracking#, WayBill)
Set-RerouteCell "V9" "51540748"
Set-RerouteCell "B9" "06-01-2022"
Set-RerouteCell "X9" '$66.51'
Set-RerouteCell "Y9" "1Z44R7R60391369101"
Set-RerouteCell "Z9" "FCUPSG1012592"

# Row 2 - new Pitney Bowes shipment (OrderId, pickUpDate, Amount, Tracking#, WayBill)
Set-RerouteCell "V2" "51540840"
Set-RerouteCell "B2" "06-02-2022"
Set-RerouteCell "X2" '$570.48'
Set-RerouteCell "Y2" "361U698177"
Set-RerouteCell "Z2" "FCPBID1001701"
